$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (2023-09-15 -> 2023-09-16) for rows 2 through 6.
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45185
}
